$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Condense / reword the paragraph describing the NodeJs backend work.
#    "...using NodeJs. These servers not only handled authentication and
#     API calls from the front end but also seamlessly interacted..."
#    becomes
#    "...using NodeJs which handled authentication and API calls from the
#     front end and also seamlessly interacted..."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "using NodeJs. These servers not only handled authentication and  API calls from the front end but also seamlessly interacted",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "using NodeJs which handled authentication and API calls from the front end and also seamlessly interacted",
    2) | Out-Null

# ---------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark (Word always records the location of the
#    most recent edit under this bookmark name) from the end of the
#    "These principles..." paragraph to the spot of the last real edit,
#    which is inside "company's" in the management-website paragraph.
#    Adding a bookmark with a name that already exists simply relocates
#    it, so the stale one at the end of the document is removed.
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("was a management website tailored for the compa", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r) | Out-Null

$d.Save()
